$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.493.15"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.89%  "
$ws.Range("D3").Value = "'1.917.22"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.14%  "
$ws.Range("E4").Value = "  +0.39%  "
$ws.Range("D5").Value = "'239.19"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.55%  "
$ws.Range("E6").Value = "  +0.25%  "
$ws.Range("D7").Value = "'0.4778"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.14%  "
$ws.Range("D8").Value = "'0.2870"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.12%  "
$ws.Range("D9").Value = "'0.06678"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.40%  "
$ws.Range("D10").Value = "'18.89"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.03%  "
$ws.Range("D11").Value = "'102.80"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.77%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "'0.07737"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.81%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "'1.919.02"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.51%  "
$ws.Range("D14").Value = "'5.223"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.68%  "
$ws.Range("D15").Value = "'0.6794"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.27%  "
$ws.Range("D16").Value = "'262.40"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -7.22%  "
$ws.Range("D17").Value = "'30.511.34"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.90%  "
$ws.Range("E18").Value = "  +0.35%  "
$ws.Range("D19").Value = "'0.000007493"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.05%  "
$ws.Range("D20").Value = "'12.72"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.23%  "
$ws.Range("D21").Value = "'5.416"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.73%  "
$ws.Range("E22").Value = "  +0.49%  "
$ws.Range("D23").Value = "'6.336"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.67%  "
$ws.Range("D24").Value = "'9.507"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.42%  "
$ws.Range("D25").Value = "'163.12"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Value = "'18.98"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.20%  "
$ws.Range("D27").Value = "'2.099"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.20%  "
$ws.Range("D28").Value = "'0.1015"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.86%  "
$ws.Range("D29").Value = "'1.384"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.28%  "
$ws.Range("D30").Value = "'4.582"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.49%  "
$ws.Range("D31").Value = "'1.516"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.23%  "
$ws.Range("D32").Value = "'4.237"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.84%  "
$ws.Range("D33").Value = "'0.04799"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.90%  "
$ws.Range("D34").Value = "'0.7327"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.90%  "
$ws.Range("D35").Value = "'1.121"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.52%  "
$ws.Range("D36").Value = "'1.002"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.18%  "
$ws.Range("D37").Value = "'2.681"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.94%  "
$ws.Range("D38").Value = "'0.01927"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.46%  "
$ws.Range("D39").Value = "'2.639"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.29%  "
$ws.Range("D40").Value = "'6.306"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.15%  "
$ws.Range("D41").Value = "'74.92"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.08%  "
$ws.Range("D42").Value = "'1.994"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.31%  "
$ws.Range("D43").Value = "'106.88"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.99%  "
$ws.Range("D44").Value = "'0.8613"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.40%  "
$ws.Range("D45").Value = "'0.4273"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.57%  "
$ws.Range("E46").Value = "  +0.18%  "
$ws.Range("D47").Value = "'1.009.20"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.64%  "
$ws.Range("D48").Value = "'7.499"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -8.63%  "
$ws.Range("D49").Value = "'0.1201"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.67%  "
$ws.Range("D50").Value = "'35.17"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.29%  "
$ws.Range("D51").Value = "'8.861"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.35%  "
